$d = $word.ActiveDocument

# Build the OOXML fragment (wrapped as a minimal package) describing the
# new Module 4/5/6 paragraphs that get appended after the existing final
# paragraph ("I'm hoping to review the difference between on-policy and
# off-policy learning...") and before the section properties.
$fragment = @'
<w:p/>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Heading1"/>
  </w:pPr>
  <w:r>
    <w:t>Module 4</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> – Q-Learning</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:ind w:firstLine="720"/>
  </w:pPr>
  <w:r>
    <w:t>This module introduced Q-learning, which is a model-free reinforcement learning algorithm where the agent learns a Q-function to estimate the expected future reward for an action in a given state. Q-learning incorporates the temporal difference error</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> that we covered with temporal difference learning. Q-learning is an off-policy method, so the optimal policy is learned independently of the behavior policy. </w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:tab/>
    <w:t xml:space="preserve">I found the cliff walk example instructive, and I enjoyed the lab assignment. It was interesting to see how the agent would fail often in the beginning, as the optimal policy was being discovered, then fail very rarely once it learned the correct policy. </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Heading1"/>
  </w:pPr>
  <w:r>
    <w:t>Module 5</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:tab/>
    <w:t xml:space="preserve">We learned about Deep Q-Networks, which extend the q-network framework by applying deep learning to estimate the Q-function. The reading assignment for the module was a paper introducing a deep q-network used to play Atari 2600 games, which learned from video input and outperformed human experts on several games. The paper also introduced experience replay, where the agent can learn from past experiences by selecting one from its “replay memory” at random. </w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:tab/>
    <w:t>I thought this paper was one of the coolest experiments I’ve read about. Creating an agent that can outperform humans on not just one game, but generalize to beat humans on multiple games, is impressive. I found it fascinating that the researchers also trained the network with the visual input from the game itself, rather than transcribing the space into something simpler.</w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Heading1"/>
  </w:pPr>
  <w:r>
    <w:lastRenderedPageBreak/>
    <w:t>Module 6 – Extending Deep Q-Networks</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:ind w:firstLine="720"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">We learned about some improvements to the basic Q-Network architecture that have been published recently, including prioritized replay and double deep q-networks. Prioritized replay deals with getting the agent to replay experiences with high learning value to speed up learning and improve performance. Double q-learning decouples the action selection and evaluation into two separate networks, thus reducing the overestimation bias that a single deep q-network tends to have. </w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:tab/>
    <w:t xml:space="preserve">I really enjoyed the reading assignment for this module – reading the published papers seemed more real than the textbook, and the experimental examples in the paper helped with my understanding of the underlying concepts. I found the prioritized experience replay especially fascinating, partly because the concept is so intuitive – it makes sense to have the agent replay experiences where it performs poorly, or expects a different outcome than what occurs. The agent should spend more time on these experiences, because they are by definition not well understood by the agent. </w:t>
  </w:r>
</w:p>
'@

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512">' +
      '<pkg:xmlData>' +
        '<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">' +
          '<Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/>' +
        '</Relationships>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' + $fragment + '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

# Collapse to an insertion point right at the end of the body content
# (i.e. right after the last run of the final paragraph, before sectPr)
# and insert the new paragraphs there.
$r = $d.Content
$r.Collapse(0)
$r.InsertXML($xml)
